$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 8118
$ws.Range("F8").Value = 85
$ws.Range("F9").Value = 7045
$ws.Range("F10").Value = 1132
$ws.Range("F11").Value = 537
$ws.Range("F12").Value = 486
$ws.Range("F14").Value = 699
$ws.Range("F15").Value = 348
$ws.Range("F18").Value = 228
$ws.Range("F19").Value = 76
$ws.Range("F21").Value = 51
$ws.Range("F22").Value = 11544
$ws.Range("F23").Value = 3
$ws.Range("F25").Value = 2240
$ws.Range("F27").Value = 3100
$ws.Range("F28").Value = 52
$ws.Range("F29").Value = 2661
$ws.Range("F31").Value = 21
$ws.Range("F32").Value = 278
$ws.Range("F35").Value = 1601
$ws.Range("F37").Value = 96
$ws.Range("F38").Value = 5780
$ws.Range("F40").Value = 1778
$ws.Range("F42").Value = 832
$ws.Range("F44").Value = 187
$ws.Range("F47").Value = 1510
$ws.Range("F48").Value = 97

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 25
$ws.Range("F10").Value = 52
$ws.Range("F20").Value = 64

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 221
$ws.Range("F3").Value = 358

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 221
$ws.Range("F4").Value = 358
$ws.Range("F7").Value = 8118
$ws.Range("F10").Value = 85
$ws.Range("F11").Value = 7045
$ws.Range("F12").Value = 7045
$ws.Range("F13").Value = 1132
$ws.Range("F14").Value = 537
$ws.Range("F15").Value = 486
$ws.Range("F16").Value = 699
$ws.Range("F17").Value = 348
$ws.Range("F20").Value = 228
$ws.Range("F22").Value = 51
$ws.Range("F25").Value = 11544
$ws.Range("F27").Value = 3
$ws.Range("F29").Value = 2240
$ws.Range("F30").Value = 2240
$ws.Range("F31").Value = 3100
$ws.Range("F32").Value = 2661
$ws.Range("F33").Value = 21
$ws.Range("F34").Value = 278
$ws.Range("F38").Value = 1601
$ws.Range("F40").Value = 96
$ws.Range("F41").Value = 5780
$ws.Range("F42").Value = 64
$ws.Range("F43").Value = 1778
$ws.Range("F46").Value = 832
$ws.Range("F47").Value = 187
$ws.Range("F50").Value = 1510
